# Update correlation results for "all_tools" and "typestate_checker" sheets
# (JaTyC working on DS9), per commit "Update corr results with JaTyC working on DS9"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: all_tools
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all_tools")

# Column widths for columns J (10) and L (12) widen slightly (to match
# columns I/K at 21.7109375 chars). The host's ColumnWidth setter snaps to
# a 1/6-character grid, so 20.8 is the closest input that lands on the
# nearest reachable width to 21.7109375.
$wsAll.Columns.Item(10).ColumnWidth = 20.8
$wsAll.Columns.Item(12).ColumnWidth = 20.8

# Row 13
$wsAll.Cells.Item(13, 7).Value = 219
$wsAll.Cells.Item(13, 9).Value = -0.4048169271682026
$wsAll.Cells.Item(13, 10).Value = 0.003684245901134024
$wsAll.Cells.Item(13, 11).Value = -0.5344276234347691
$wsAll.Cells.Item(13, 12).Value = 0.002348131666704992

# Row 14
$wsAll.Cells.Item(14, 7).Value = 219
$wsAll.Cells.Item(14, 9).Value = 0.1515153313763357
$wsAll.Cells.Item(14, 10).Value = 0.2780537029565812
$wsAll.Cells.Item(14, 11).Value = 0.2312941947390607
$wsAll.Cells.Item(14, 12).Value = 0.2187797748162401

# Row 15
$wsAll.Cells.Item(15, 7).Value = 219
$wsAll.Cells.Item(15, 9).Value = 0.2373626514505708
$wsAll.Cells.Item(15, 10).Value = 0.09051426625460415
$wsAll.Cells.Item(15, 11).Value = 0.3315718656719684
$wsAll.Cells.Item(15, 12).Value = 0.07346214270566978

# Row 16
$wsAll.Cells.Item(16, 7).Value = 219
$wsAll.Cells.Item(16, 9).Value = -0.3966737953014234
$wsAll.Cells.Item(16, 10).Value = 0.004392923006413007
$wsAll.Cells.Item(16, 11).Value = -0.5638269314967141
$wsAll.Cells.Item(16, 12).Value = 0.001175300429566541

# ---------------------------------------------------------------------------
# Sheet: typestate_checker
# ---------------------------------------------------------------------------
$wsTs = $wb.Worksheets.Item("typestate_checker")

# Row 13
$wsTs.Cells.Item(13, 6).Value = 30
$wsTs.Cells.Item(13, 7).Value = 111
$wsTs.Cells.Item(13, 9).Value = -0.3137312989174262
$wsTs.Cells.Item(13, 10).Value = 0.02848639994432931
$wsTs.Cells.Item(13, 11).Value = -0.4131284275634866
$wsTs.Cells.Item(13, 12).Value = 0.02326314515800765

# Row 14
$wsTs.Cells.Item(14, 6).Value = 30
$wsTs.Cells.Item(14, 7).Value = 111
$wsTs.Cells.Item(14, 9).Value = 0.14638501094228
$wsTs.Cells.Item(14, 10).Value = 0.3077335885571051
$wsTs.Cells.Item(14, 11).Value = 0.2212047859710532
$wsTs.Cells.Item(14, 12).Value = 0.2401033542701417

# Row 15
$wsTs.Cells.Item(15, 6).Value = 30
$wsTs.Cells.Item(15, 7).Value = 111
$wsTs.Cells.Item(15, 9).Value = 0.2124307787987451
$wsTs.Cells.Item(15, 10).Value = 0.1403732278087201
$wsTs.Cells.Item(15, 11).Value = 0.2913831605167395
$wsTs.Cells.Item(15, 12).Value = 0.1182228233126965

# Row 16
$wsTs.Cells.Item(16, 6).Value = 30
$wsTs.Cells.Item(16, 7).Value = 111
$wsTs.Cells.Item(16, 9).Value = -0.3809072222823607
$wsTs.Cells.Item(16, 10).Value = 0.007761477051653751
$wsTs.Cells.Item(16, 11).Value = -0.5322611581053379
$wsTs.Cells.Item(16, 12).Value = 0.002464998415486295
